$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("A4").Value = 1

$ws.Range("B4:G4").NumberFormat = "@"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("B4").Value = "FARCOVIT B12 30 CAPS."
$ws.Range("N4").Value = "0:0"

$ws.Range("H4:K4").NumberFormat = "@"
$ws.Range("H4").Value = "9:1"

$ws.Range("L4").Value = 15

# Row 5 updates
$ws.Rows("5").RowHeight = 26.25
$ws.Range("K5").Value = 15
